$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Text content changes (shared-string edits)
# ---------------------------------------------------------------
$ws.Cells.Item(4, 1).Value  = "АБК ЖДЦ"
$ws.Cells.Item(17, 1).Value = "Травнсформатор власних потреб (ТМ-40)"
$ws.Cells.Item(18, 1).Value = "Їдальня (ТМ-40)"
$ws.Cells.Item(19, 1).Value = "Фікальна насосна (ТМ-160)"

# Power value for the "Фікальна насосна" transformer changed 120 -> 160 kVA
$ws.Range("C19").Value = 160

# ---------------------------------------------------------------
# 2. Column A gets narrower, rows wrap the (now longer) names
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 17

# Apply the existing "wrap text only" format (same as H2/I2) to column A
# of the data rows, without touching the shared default style of other
# sheets/cells.
$ws.Range("H2").Copy()
$ws.Range("A3:A15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A17:A19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row heights for the rows whose wrapped text now needs two lines
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 45
$ws.Rows.Item(19).RowHeight = 30

# A16 ("Споживач") gets the centred + wrapped header look
$ws.Range("B16").Copy()
$ws.Range("A16").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A16").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A16").VerticalAlignment = -4108    # xlCenter
$ws.Range("A16").WrapText = $true

# ---------------------------------------------------------------
# 3. View state: scrolled further down, new selected cell
# ---------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("D23").Select()

Write-Host "edit complete"
